$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "Payload Budget" sheet between "Mass budget" and
#    "Component weights".
# ---------------------------------------------------------------------------
$massBudget = $wb.Worksheets.Item("Mass budget")
$ws = $wb.Worksheets.Add($null, $massBudget)
$ws.Name = "Payload Budget"

# ---------------------------------------------------------------------------
# 2. Layout / column widths.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 25.6640625

# ---------------------------------------------------------------------------
# 3. Header row.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Component "
$ws.Range("B2").Value = "Mass (g)"
$ws.Range("C2").Value = "Cost (EUR)"
$ws.Range("A2:C2").Font.Name = "Arial"
$ws.Range("A2:C2").Font.Size = 12
$ws.Range("A2:C2").Font.Bold = $true
$ws.Range("A2:C2").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 4. Data rows (component / mass / cost).
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Autopilot"
$ws.Range("B3").Value = 10.54
$ws.Range("C3").Value = 100

$ws.Range("A4").Value = "Pitot tube"
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 20

$ws.Range("A5").Value = "Servomotors x 3"
$ws.Range("B5").Formula = "=4.5*3"
$ws.Range("C5").Value = 90

$ws.Range("A6").Value = "UBLOX"
$ws.Range("B6").Value = 23
$ws.Range("C6").Value = 350
$ws.Range("D6").Value = "https://drotek.com/shop/en/drotek-parts/792-xl-rtk-gps-neo-m8p-rover.html"

$ws.Range("A7").Value = "Xbee"
$ws.Range("B7").Value = 4
$ws.Range("C7").Value = 50
$ws.Range("D7").Value = "http://www.mouser.ch/ProductDetail/Digi-International/XB8-DMUS-002/?qs=%2fha2pyFaduhamZ1j%2fTvqvLufopfOt%252bEr7wrfpr46JtwlpBsVij4AWA%3d%3d"

$ws.Range("A8").Value = "Battery"
$ws.Range("B8").Value = 36
$ws.Range("C8").Value = 5

$ws.Range("A9").Value = "Cabling"
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 2

$ws.Range("A10").Value = "Mounting"
$ws.Range("B10").Value = 10
$ws.Range("C10").Value = 3

$ws.Range("A11").Value = "Antenna"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 10

$ws.Range("A12").Value = "Wing"
$ws.Range("B12").Value = 100
$ws.Range("C12").Value = 30

$ws.Range("A13").Value = "Hinges + locking mech"
$ws.Range("B13").Value = 10
$ws.Range("C13").Value = 4

$ws.Range("A14").Value = "Camera"
$ws.Range("B14").Value = 20
$ws.Range("C14").Value = 30

# Give the whole data block the Arial 12 look.
$ws.Range("A3:C14").Font.Name = "Arial"
$ws.Range("A3:C14").Font.Size = 12

# Highlight the UBLOX / Xbee rows in green, like the source sheet.
$ws.Range("A6:C7").Interior.Color = 65280

# ---------------------------------------------------------------------------
# 5. Totals.
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "Total without margin"
$ws.Range("B15:C15").Formula = "=SUM(B3:B14)"

$ws.Range("A16").Value = "Total with 20% margin"
$ws.Range("B16:C16").Formula = "=B15+20/100*B15"

$ws.Range("A15:C16").Font.Name = "Arial"
$ws.Range("A15:C16").Font.Size = 12
$ws.Range("A15:C16").Font.Bold = $true

# ---------------------------------------------------------------------------
# 6. Ground station line, further down the sheet.
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = "Ground Station GPS RTK"
$ws.Range("C21").Value = 350
$ws.Range("A21").Font.Name = "Arial"
$ws.Range("A21").Font.Size = 12
$ws.Range("C21").Font.Name = "Arial"
$ws.Range("C21").Font.Size = 12

# ---------------------------------------------------------------------------
# 7. Hyperlinks for the two part-number references.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("D6"), "https://drotek.com/shop/en/drotek-parts/792-xl-rtk-gps-neo-m8p-rover.html")
$ws.Hyperlinks.Add($ws.Range("D7"), "http://www.mouser.ch/ProductDetail/Digi-International/XB8-DMUS-002/?qs=%2fha2pyFaduhamZ1j%2fTvqvLufopfOt%252bEr7wrfpr46JtwlpBsVij4AWA%3d%3d")

# D7's link keeps the small underlined-blue Arial look used elsewhere in the
# workbook rather than the default (larger) Hyperlink style.
$ws.Range("D7").Font.Name = "Arial"
$ws.Range("D7").Font.Size = 10
$ws.Range("D7").Font.Underline = $true
$ws.Range("D7").Font.Color = 16711680

# ---------------------------------------------------------------------------
# 8. Selection / activation state, matching the authored file.
# ---------------------------------------------------------------------------
$ws.Range("D17").Select()
$ws.Activate()
